$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1073.7894
$ws.Range("I107").Value = 500.07693
$ws.Range("J107").Value = 2316.8333
$ws.Range("K107").Value = 500.07693
$ws.Range("L107").Value = 2316.8333
$ws.Range("M107").Value = 1419.92307
$ws.Range("N107").Value = -6156.8333
$ws.Range("H112").Value = 1884.7255
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 1900.4286
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 5701.2858
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -7917.2858
$ws.Range("H116").Value = 12000.9
$ws.Range("I116").Value = 26874.75
$ws.Range("J116").Value = 2085
$ws.Range("K116").Value = 26874.75
$ws.Range("L116").Value = 2085
$ws.Range("M116").Value = -23432.75
$ws.Range("N116").Value = -8969
$ws.Range("H121").Value = 738.2
$ws.Range("I121").Value = 199.5
$ws.Range("J121").Value = 872.875
$ws.Range("K121").Value = 598.5
$ws.Range("L121").Value = 2618.625
$ws.Range("M121").Value = 1148.5
$ws.Range("N121").Value = -6112.625
$ws.Range("H132").Value = 863.97675
$ws.Range("I132").Value = 803.0513
$ws.Range("K132").Value = 2409.1539
$ws.Range("M132").Value = 120.8461000000002
$ws.Range("H137").Value = 1174.375
$ws.Range("I137").Value = 715.4167
$ws.Range("K137").Value = 2146.2501
$ws.Range("M137").Value = 403.7498999999998
$ws.Range("H138").Value = 2408.0703
$ws.Range("I138").Value = 3085.1155
$ws.Range("J138").Value = 2016.8889
$ws.Range("K138").Value = 9255.3465
$ws.Range("L138").Value = 6050.6667
$ws.Range("M138").Value = -4115.3465
$ws.Range("N138").Value = -16330.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 505616.2
$ws.Range("I2").Value = 695008.1
$ws.Range("J2").Value = 571
$ws.Range("K2").Value = 695008.1
$ws.Range("L2").Value = 571
$ws.Range("M2").Value = -694895.1
$ws.Range("N2").Value = -797
$ws.Range("H16").Value = 1849.75
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1849.75
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1849.75
$ws.Range("N16").Value = -2423.75
$ws.Range("H32").Value = 5860.3877
$ws.Range("I32").Value = 4469.07
$ws.Range("K32").Value = 4469.07
$ws.Range("M32").Value = -4182.07
$ws.Range("H61").Value = 6599.885
$ws.Range("I61").Value = 6802.353
$ws.Range("J61").Value = 6217.4443
$ws.Range("K61").Value = 6802.353
$ws.Range("L61").Value = 6217.4443
$ws.Range("M61").Value = -6590.353
$ws.Range("N61").Value = -6641.4443
$ws.Range("H74").Value = 829.6429000000001
$ws.Range("I74").Value = 549.91895
$ws.Range("K74").Value = 549.91895
$ws.Range("M74").Value = 324.08105
$ws.Range("H77").Value = 829.6429000000001
$ws.Range("I77").Value = 549.91895
$ws.Range("K77").Value = 2749.59475
$ws.Range("M77").Value = 1618.40525
$ws.Range("H110").Value = 478
$ws.Range("I110").Value = 478
$ws.Range("K110").Value = 478
$ws.Range("M110").Value = 1567
$ws.Range("H116").Value = 505616.2
$ws.Range("I116").Value = 695008.1
$ws.Range("J116").Value = 571
$ws.Range("K116").Value = 695008.1
$ws.Range("L116").Value = 571
$ws.Range("M116").Value = -692714.1
$ws.Range("N116").Value = -5159
$ws.Range("H136").Value = 6599.885
$ws.Range("I136").Value = 6802.353
$ws.Range("J136").Value = 6217.4443
$ws.Range("K136").Value = 20407.059
$ws.Range("L136").Value = 18652.3329
$ws.Range("M136").Value = -17857.059
$ws.Range("N136").Value = -23752.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 505616.2
$ws.Range("I3").Value = 695008.1
$ws.Range("J3").Value = 571
$ws.Range("K3").Value = 695008.1
$ws.Range("L3").Value = 571
$ws.Range("M3").Value = -694894.1
$ws.Range("N3").Value = -799
$ws.Range("H20").Value = 2250.2
$ws.Range("I20").Value = 1910.875
$ws.Range("K20").Value = 1910.875
$ws.Range("M20").Value = -1663.875
$ws.Range("H81").Value = 37498.5
$ws.Range("J81").Value = 37498.5
$ws.Range("L81").Value = 37498.5
$ws.Range("N81").Value = -39620.5
$ws.Range("H84").Value = 37498.5
$ws.Range("J84").Value = 37498.5
$ws.Range("L84").Value = 112495.5
$ws.Range("N84").Value = -123103.5
$ws.Range("H105").Value = 2399.7083
$ws.Range("I105").Value = 2259.8
$ws.Range("J105").Value = 3099.25
$ws.Range("K105").Value = 2259.8
$ws.Range("L105").Value = 3099.25
$ws.Range("M105").Value = -512.8000000000002
$ws.Range("N105").Value = -6593.25
$ws.Range("H135").Value = 57427.57
$ws.Range("J135").Value = 57427.57
$ws.Range("L135").Value = 57427.57
$ws.Range("N135").Value = -67567.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2753.4092
$ws.Range("I31").Value = 2420.3333
$ws.Range("K31").Value = 2420.3333
$ws.Range("M31").Value = -2125.3333
$ws.Range("H34").Value = 2753.4092
$ws.Range("I34").Value = 2420.3333
$ws.Range("K34").Value = 2420.3333
$ws.Range("M34").Value = -2218.3333
$ws.Range("H134").Value = 2117.2188
$ws.Range("I134").Value = 1981.0385
$ws.Range("K134").Value = 5943.1155
$ws.Range("M134").Value = -3408.1155
$ws.Range("H141").Value = 59284
$ws.Range("J141").Value = 57164.668
$ws.Range("L141").Value = 57164.668
$ws.Range("N141").Value = -67524.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 648.41174
$ws.Range("I5").Value = 562.25
$ws.Range("K5").Value = 1686.75
$ws.Range("M5").Value = -1574.75
$ws.Range("H33").Value = 129.11111
$ws.Range("I33").Value = 133
$ws.Range("J33").Value = 121.333336
$ws.Range("K33").Value = 798
$ws.Range("L33").Value = 728.000016
$ws.Range("M33").Value = -515
$ws.Range("N33").Value = -1294.000016
$ws.Range("H118").Value = 1314.2
$ws.Range("I118").Value = 631.4
$ws.Range("J118").Value = 1997
$ws.Range("K118").Value = 1894.2
$ws.Range("L118").Value = 5991
$ws.Range("M118").Value = -651.1999999999998
$ws.Range("N118").Value = -8477
$ws.Range("H131").Value = 15328.2295
$ws.Range("J131").Value = 15967.935
$ws.Range("L131").Value = 47903.805
$ws.Range("N131").Value = -57983.805
$ws.Range("H135").Value = 648.41174
$ws.Range("I135").Value = 562.25
$ws.Range("K135").Value = 5060.25
$ws.Range("M135").Value = -2525.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3294.5
$ws.Range("I7").Value = 3439.5
$ws.Range("K7").Value = 3439.5
$ws.Range("M7").Value = -3327.5
$ws.Range("H46").Value = 1866.625
$ws.Range("I46").Value = 1198.375
$ws.Range("J46").Value = 2534.875
$ws.Range("K46").Value = 1198.375
$ws.Range("L46").Value = 2534.875
$ws.Range("M46").Value = -1010.375
$ws.Range("N46").Value = -2910.875
$ws.Range("H61").Value = 2562.125
$ws.Range("I61").Value = 2335.5908
$ws.Range("J61").Value = 3060.5
$ws.Range("K61").Value = 2335.5908
$ws.Range("L61").Value = 3060.5
$ws.Range("M61").Value = -2133.5908
$ws.Range("N61").Value = -3464.5
$ws.Range("H113").Value = 2562.125
$ws.Range("I113").Value = 2335.5908
$ws.Range("J113").Value = 3060.5
$ws.Range("K113").Value = 2335.5908
$ws.Range("L113").Value = 3060.5
$ws.Range("M113").Value = -165.5907999999999
$ws.Range("N113").Value = -7400.5
$ws.Range("H126").Value = 3294.5
$ws.Range("I126").Value = 3439.5
$ws.Range("K126").Value = 10318.5
$ws.Range("M126").Value = -7848.5
$ws.Range("H132").Value = 2406.122
$ws.Range("I132").Value = 1278.5
$ws.Range("K132").Value = 3835.5
$ws.Range("M132").Value = -1305.5
$ws.Range("H136").Value = 1611.6
$ws.Range("I136").Value = 1111.7273
$ws.Range("J136").Value = 2986.25
$ws.Range("K136").Value = 3335.1819
$ws.Range("L136").Value = 8958.75
$ws.Range("M136").Value = -785.1819
$ws.Range("N136").Value = -14058.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
$ws.Range("H107").Value = 653.7857
$ws.Range("J107").Value = 916.9
$ws.Range("L107").Value = 2750.7
$ws.Range("N107").Value = -6590.7
$ws.Range("H122").Value = 58158.5
$ws.Range("I122").Value = 132633.33
$ws.Range("J122").Value = 2302.375
$ws.Range("K122").Value = 397899.99
$ws.Range("L122").Value = 6907.125
$ws.Range("M122").Value = -395449.99
$ws.Range("N122").Value = -11807.125
$ws.Range("H136").Value = 14621054
$ws.Range("I136").Value = 19157918
$ws.Range("K136").Value = 57473754
$ws.Range("M136").Value = -57471204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M16").ClearContents()